# Update specific hashcode values in column B (column A has the code identifiers)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=9; Code="05-050305TC"; Hash="e93effb58e5970f605ae07ea0fd6480b"}
    @{Row=17; Code="05-050305TP"; Hash="bb451ec4926ef9a76c82b3a70560c0a5"}
    @{Row=34; Code="05-050316TP"; Hash="0d3cbd5cf9a3bf3ff616ce16adc4567b"}
    @{Row=79; Code="05-050307TP"; Hash="f21b89ca06c5305e282a3da162a1ba2f"}
    @{Row=126; Code="05-050309A"; Hash="30992a194a56e3775d7bc9fa5a64bc24"}
    @{Row=136; Code="05-050312TC"; Hash="5e3fe43d9be5b777179b6c69eea2d63f"}
    @{Row=150; Code="05-050307A"; Hash="a7846e86e6fbfe8b5adf41eb2691103c"}
    @{Row=159; Code="05-050203TP"; Hash="17e6f09fd8ea8a8972bc475df817080f"}
    @{Row=162; Code="05-050308A"; Hash="e0e40dd369b501d7a760405fc16826ff"}
    @{Row=169; Code="05-050203TC"; Hash="6afcb86346c0f16cac73003425cae14d"}
    @{Row=175; Code="05-050303TP"; Hash="0e80f321852f84a5589ee4ed776865b7"}
    @{Row=180; Code="05-050303TC"; Hash="2851dfb7ee2096dfbfb00ec7833f6902"}
    @{Row=183; Code="05-050305A"; Hash="0eb561f33a5f6140b3aa69014266153b"}
    @{Row=200; Code="05-050306A"; Hash="f80e0ff992d99c2f1ced4e5b5e305976"}
    @{Row=213; Code="05-050303A"; Hash="247ce9bcf77ab8b1ad96f59718f44de5"}
    @{Row=228; Code="05-050304A"; Hash="d0aec2c73e440a1866e6bbb0730408de"}
    @{Row=281; Code="05-050201TC"; Hash="91d6cecafdef3ad37838abc58fd1f3c8"}
    @{Row=302; Code="05-050310TP"; Hash="d263c9cd625e0cc36308d3fec4350e23"}
    @{Row=339; Code="05-050201TP"; Hash="1e506b1f2a033ed20095cbdd53afc20a"}
    @{Row=460; Code="05-050204A"; Hash="0cd8625297c32aba25b0f61545f1b53e"}
    @{Row=461; Code="05-050313A"; Hash="1987c093b7249b83a32c4f426a98e594"}
    @{Row=500; Code="05-050202A"; Hash="59328d6fbee2ac587678815c09af1874"}
    @{Row=501; Code="05-050311A"; Hash="2f3dfc70d7f041da9765e62f76ca913a"}
    @{Row=506; Code="05-050306TP"; Hash="857f16fc79e5fc5ffa0511f91f7b30ce"}
    @{Row=514; Code="05-050317TC"; Hash="1522a941e7773172e4dd4ad354ab0470"}
    @{Row=517; Code="05-050203A"; Hash="4411e56c2ff7e6ec8787d8f6be166e8b"}
    @{Row=524; Code="05-050317TP"; Hash="929b51ea954a9711462847af84dc8432"}
    @{Row=547; Code="05-050201A"; Hash="f6e4456a75f8fa1f2b7ad0cc3469b942"}
    @{Row=550; Code="05-050310A"; Hash="345984d1f1a72d556b2fb2538b0e94aa"}
    @{Row=563; Code="05-050308TC"; Hash="c110054283d2d57b80f1cdba7cc6ce42"}
    @{Row=572; Code="05-050308TP"; Hash="2829c5fc1f67e224165dc8d654e289f4"}
    @{Row=616; Code="05-050204TP"; Hash="cf51451dd6f5b3073cd680b0a9c8f098"}
    @{Row=627; Code="05-050204TC"; Hash="cd0f810a0814b71df06adc86d49f9165"}
    @{Row=666; Code="05-050317A"; Hash="abf90ea370bd45b91b48fbc900bc506d"}
    @{Row=700; Code="05-050304TC"; Hash="c1be0d083ce0ad19eb1f14e63dd5771f"}
    @{Row=715; Code="05-050304TP"; Hash="d6ec5b2a28c05cafb949242c8f5515d0"}
    @{Row=729; Code="05-050316A"; Hash="52d45121b8d9764e0fdb39e8ce4c0c5e"}
    @{Row=819; Code="05-050202TP"; Hash="f918429f8f38492013789bfd11f54108"}
    @{Row=830; Code="05-050311TC"; Hash="39131b3cfdad3487567b097fc174ea20"}
    @{Row=835; Code="05-050311TP"; Hash="6c0c01f6b02ef111a430a37b418b5556"}
    @{Row=854; Code="05-050309TC"; Hash="78f94801636925cb39037cc331070e59"}
    @{Row=862; Code="05-050309TP"; Hash="15adcc8626573003a2667afe259f8d2e"}
)

foreach ($u in $updates) {
    $actualCode = $ws.Cells.Item($u.Row, 1).Value2
    if ($actualCode -ne $u.Code) {
        Write-Output "WARNING: row $($u.Row) expected code '$($u.Code)' but found '$actualCode'"
    }
    $ws.Cells.Item($u.Row, 2).Value = $u.Hash
}

